# Demo.xlsx edit: turn the 3-sheet scratch workbook into a "NPC" table
# definition workbook (NPC data sheet + @Types metadata sheet), matching
# a typical Luban/table-export config layout.

$wb = $excel.ActiveWorkbook

# --- Sheet1 -> "NPC" --------------------------------------------------
$npc = $wb.Worksheets.Item(1)
$npc.Name = "NPC"

# Header / type / tag / comment rows (row 1-4), then one data row (row 5)
$npc.Range("A1").Value = "ID"
$npc.Range("B1").Value = "Name"
$npc.Range("C1").Value = "HP"
$npc.Range("D1").Value = "Skills"

$npc.Range("A2").Value = "int32"
$npc.Range("B2").Value = "string"
$npc.Range("C2").Value = "float"
$npc.Range("D2").Value = "repeated int32"

$npc.Range("A3").Value = "RepeatCheck:true"
$npc.Range("B3").Value = "MakeIndex:true"
$npc.Range("D3").Value = 'ListSpliter:","'

$npc.Range("A4").Value = "ID"
$npc.Range("B4").Value = "名称"
$npc.Range("C4").Value = "血量上限"
$npc.Range("D4").Value = "技能列表"

$npc.Range("A5").Value = 10001
$npc.Range("B5").Value = "npc1"
$npc.Range("C5").Value = 1000
$npc.Range("D5").Value = "1,2,3,4"

# Column widths (character units). The host quantizes stored width to
# 1/7-character steps, so these are the nearest achievable inputs to the
# authored widths of 23.125 / 20.5 / 24.75 / 19.
$npc.Columns.Item(1).ColumnWidth = 22.428571428571427
$npc.Columns.Item(2).ColumnWidth = 19.857142857142858
$npc.Columns.Item(3).ColumnWidth = 24.0
$npc.Columns.Item(4).ColumnWidth = 18.285714285714285

# --- Sheet2 -> "@Types" -------------------------------------------------
$types = $wb.Worksheets.Item(2)
$types.Name = "@Types"
$types.Range("A1").Value = 'TableName: "NPC" Package: "table" CSClassHeader: "[System.Serializable]"'

# --- Selection / active tab state --------------------------------------
# Final diff state: NPC sheet not active, selection parked at A5;
# @Types sheet is the active tab with selection at A2.
[void]$npc.Range("A5").Select()
[void]$types.Activate()
[void]$types.Range("A2").Select()
